$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The three brand-new shared strings ("50", "KK", "Feb-17") must be created
# in this exact order so they land at shared-string indices 206, 207, 208
# (matching the target workbook). We touch the cells that introduce them
# first, then fill in the remaining (already-existing) string values.
# ---------------------------------------------------------------------------

# J78 = "50"  (new shared string -> 206)
$ws.Cells.Item(78, 10).Value = "50"
$ws.Cells.Item(78, 10).NumberFormat = "@"
$ws.Cells.Item(78, 10).Font.Size = 12
$ws.Cells.Item(78, 10).HorizontalAlignment = -4152

# K77 = "KK"  (new shared string -> 207)
$ws.Cells.Item(77, 11).Value = "KK"
$ws.Cells.Item(77, 11).Font.Size = 16

# N15 = "Feb-17"  (new shared string -> 208)
$ws.Cells.Item(15, 14).Value = "Feb-17"
$ws.Cells.Item(15, 14).NumberFormat = "@"
$ws.Cells.Item(15, 14).Font.Size = 16

# ---------------------------------------------------------------------------
# Mark cases on rows 15, 16, 32, 49, 66 as Discharged on Feb-17.
# ---------------------------------------------------------------------------

$ws.Cells.Item(15, 11).Value = "Discharged"

$ws.Cells.Item(16, 11).Value = "Discharged"
$ws.Cells.Item(16, 14).Value = "Feb-17"
$ws.Cells.Item(16, 14).NumberFormat = "@"
$ws.Cells.Item(16, 14).Font.Size = 16

$ws.Cells.Item(32, 11).Value = "Discharged"
$ws.Cells.Item(32, 14).Value = "Feb-17"
$ws.Cells.Item(32, 14).NumberFormat = "@"
$ws.Cells.Item(32, 14).Font.Size = 16

$ws.Cells.Item(49, 11).Value = "Discharged"
$ws.Cells.Item(49, 14).Value = "Feb-17"
$ws.Cells.Item(49, 14).NumberFormat = "@"
$ws.Cells.Item(49, 14).Font.Size = 16

$ws.Cells.Item(66, 11).Value = "Discharged"
$ws.Cells.Item(66, 14).Value = "Feb-17"
$ws.Cells.Item(66, 14).NumberFormat = "@"
$ws.Cells.Item(66, 14).Font.Size = 16

# ---------------------------------------------------------------------------
# Append two new case rows (77 and 78).
# ---------------------------------------------------------------------------

$ws.Rows.Item(77).RowHeight = 21
$ws.Rows.Item(78).RowHeight = 21

# Row 77
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 1).Font.Size = 16

$ws.Cells.Item(77, 2).Value = 1.310438
$ws.Cells.Item(77, 2).Font.Size = 16

$ws.Cells.Item(77, 3).Value = 103.84718100000001

$ws.Cells.Item(77, 4).Value = "Feb-16"
$ws.Cells.Item(77, 4).NumberFormat = "@"
$ws.Cells.Item(77, 4).HorizontalAlignment = -4152

$ws.Cells.Item(77, 5).Value = 1
$ws.Cells.Item(77, 5).Font.Size = 16

$ws.Cells.Item(77, 6).Value = "Male"
$ws.Cells.Item(77, 6).Font.Size = 16

$ws.Cells.Item(77, 7).Value = "China, Wuhan"
$ws.Cells.Item(77, 7).Font.Size = 16

$ws.Cells.Item(77, 12).Value = "Singaporean"
$ws.Cells.Item(77, 12).Font.Size = 16

$ws.Cells.Item(77, 13).Value = "Feb-09"
$ws.Cells.Item(77, 13).NumberFormat = "@"
$ws.Cells.Item(77, 13).Font.Size = 12

# Row 78
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 1).Font.Size = 16

$ws.Cells.Item(78, 2).Value = 1.3535889999999999
$ws.Cells.Item(78, 2).Font.Size = 16

$ws.Cells.Item(78, 3).Value = 103.859931

$ws.Cells.Item(78, 4).Value = "Feb-17"
$ws.Cells.Item(78, 4).NumberFormat = "@"
$ws.Cells.Item(78, 4).HorizontalAlignment = -4152

$ws.Cells.Item(78, 5).Value = 35
$ws.Cells.Item(78, 5).Font.Size = 16

$ws.Cells.Item(78, 6).Value = "Male"
$ws.Cells.Item(78, 6).Font.Size = 16

$ws.Cells.Item(78, 7).Value = "Singapore"
$ws.Cells.Item(78, 7).Font.Size = 16

$ws.Cells.Item(78, 11).Value = "NCID"
$ws.Cells.Item(78, 11).Font.Size = 16

$ws.Cells.Item(78, 12).Value = "Singaporean"
$ws.Cells.Item(78, 12).Font.Size = 16

# ---------------------------------------------------------------------------
# Update the view so it matches the scrolled/selected state in the target.
# ---------------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 67
$ws.Range("C78").Select() | Out-Null
